$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All source cells are plain text (t="inlineStr"), even the ones that look
# like plain numbers (prices such as "576.73" or "0.0359"). Assigning those
# strings straight to .Value lets Excel auto-detect them as numbers, which
# both changes the cell's stored type and can silently drop formatting such
# as a trailing zero ("384.10" -> 384.1). To keep them as text exactly like
# the original file, force the Text number format before writing the value,
# then restore the default "Normal" style so no extra formatting is left
# behind on cells that did not carry one originally.

$ws.Range('D2').Value = '66.167.25'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '3.026.09'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.66'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.28%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.027.05'
$ws.Range('E8').Value = '  +0.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.521'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('E11').Value = '  -2.01%  '
$ws.Range('E12').Value = '  +4.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000248'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.00%  '
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '66.228.54'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').Value = '3.524.28'
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.26'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.32%  '
$ws.Range('D19').Value = '3.025.95'
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '466.76'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.709'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.51'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.63'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.44%  '
$ws.Range('E25').Value = '  +0.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.47'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.46'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.49%  '
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range('E32').Value = '  -3.94%  '
$ws.Range('E33').Value = '  +6.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.28'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.42%  '
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.87'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('B37').Value = 'Mantle'
$ws.Range('C37').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.992'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.19'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +10.63%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.314'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.64'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.04'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.19%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.122'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.89'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.67'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.38%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0359'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '384.10'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.88%  '
$ws.Range('D47').Value = '2.727.05'
$ws.Range('E47').Value = '  -2.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.64'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.80'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.44%  '
$ws.Range('E51').Value = '  +3.97%  '
